# Fix Training Data Issue: the BF column ("Date") held a malformed date
# string "4-8-2012-13" (the season folder name leaking into the date
# field). NBA stats for that game were actually from 2013-04-08, so
# rewrite the column to the correct ISO date string.
#
# The replacement must stay a literal text string (not get silently
# reinterpreted as a date serial number by Excel's smart-entry), so the
# new value is entered with a leading apostrophe - the normal Excel
# "force text" convention - which keeps Value2 / the stored cell content
# exactly "2013-04-08".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "4-8-2012-13"
$newValue = "2013-04-08"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$dateCol = 58  # Column BF

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = "'" + $newValue
    }
}
